$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 504 ("「眠りたいです」" post), shifting all rows below up by one.
$ws.Rows.Item(504).Delete()
